$wb = $excel.ActiveWorkbook

# Rename the last sheet from "Sheet1" to "Insert Values"
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Insert Values"

# Make it the active sheet and set the active cell selection to Q25
$ws.Activate()
$ws.Range("Q25").Select()
